# The presentation's slide-master theme ("Integral") is re-coloured to the
# stock Office default colour palette (the author swapped the slide-master
# theme onto the plain "Office Theme" colours while leaving the font/format
# schemes untouched, since both were already the default "Office" ones).
#
# PowerPoint exposes the twelve theme colour slots through
# Design.SlideMaster.Theme.ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), each settable via .RGB (a standard VBA RGB() BGR-packed
# integer). dk1/lt1 are already black/white in both palettes, so only the
# remaining ten slots need to move to the Office defaults.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colors = $master.Theme.ThemeColorScheme

function BGR($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$colors.Item(3).RGB  = BGR 0x44 0x54 0x6A   # dk2      -> 44546A
$colors.Item(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$colors.Item(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$colors.Item(6).RGB  = BGR 0xED 0x7D 0x31   # accent2  -> ED7D31
$colors.Item(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$colors.Item(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4  -> FFC000
$colors.Item(9).RGB  = BGR 0x44 0x72 0xC4   # accent5  -> 4472C4
$colors.Item(10).RGB = BGR 0x70 0xAD 0x47   # accent6  -> 70AD47
$colors.Item(11).RGB = BGR 0x05 0x63 0xC1   # hlink    -> 0563C1
$colors.Item(12).RGB = BGR 0x95 0x4F 0x72   # folHlink -> 954F72
